# [FIX] adapt data to 8.0 compatibility
# Rename/adapt the res.partner import sample data to be compatible with
# Odoo 8.0 field values (company_type -> is_company, and the selection
# values for res.partner "type" field using their technical keys).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("res.partner")

# Header: company_type -> is_company
$ws.Range("C1").Value = "is_company"

# company_type ("Company"/"Individual") -> is_company ("yes"/"no")
$ws.Range("C2").Value = "yes"
$ws.Range("C3").Value = "no"
$ws.Range("C4").Value = "no"
$ws.Range("C5").Value = "no"
$ws.Range("C6").Value = "no"
$ws.Range("C7").Value = "yes"
$ws.Range("C8").Value = "yes"

# type values -> lowercase technical selection keys
$ws.Range("E3").Value = "contact"
$ws.Range("E4").Value = "other"
$ws.Range("E5").Value = "delivery"
$ws.Range("E6").Value = "invoice"

# Update the active selection left behind on the sheet
$ws.Range("E5").Select()
